# Increment every date in column C (rows 2-407, "Förändrad") by one day.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 407
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
